$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "23×84="; New = "49×98=" },
    @{ Old = "62×82="; New = "55×93=" },
    @{ Old = "85×13="; New = "78×73=" },
    @{ Old = "71×88="; New = "40×67=" },
    @{ Old = "47×81="; New = "44×42=" },
    @{ Old = "35×56="; New = "74×90=" },
    @{ Old = "25×61="; New = "41×45=" },
    @{ Old = "52×44="; New = "56×29=" },
    @{ Old = "19×19="; New = "73×60=" },
    @{ Old = "83×74="; New = "26×72=" },
    @{ Old = "81×49="; New = "30×90=" },
    @{ Old = "34×66="; New = "49×93=" },
    @{ Old = "97×58="; New = "77×30=" },
    @{ Old = "33×62="; New = "98×66=" },
    @{ Old = "20×77="; New = "17×38=" },
    @{ Old = "73×41="; New = "90×83=" },
    @{ Old = "63×64="; New = "52×20=" },
    @{ Old = "92×20="; New = "27×30=" },
    @{ Old = "15×55="; New = "27×57=" },
    @{ Old = "73×51="; New = "90×30=" },
    @{ Old = "43×47="; New = "57×12=" },
    @{ Old = "69×16="; New = "86×57=" },
    @{ Old = "46×26="; New = "69×26=" },
    @{ Old = "26×45="; New = "17×43=" },
    @{ Old = "61×86="; New = "36×23=" }
)

foreach ($rep in $replacements) {
    $range = $d.Content
    $range.Find.Execute($rep.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $rep.New, 2)
}
